# Iceland Premier League - re-sync match rows that were out of order.
# For a handful of same-date fixtures the match-result / odds columns
# (B, F:AC) had been written into the wrong row. This restores the
# correct row <-> data association by rotating the B,F:AC payload among
# the affected rows while leaving id (A), Div/Div Original Name (C/D)
# and Date (E) untouched (they are already correct / shared across the
# group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding the "payload" that needs to move between rows.
# (C, D, E - Div / Div Original Name / Date - stay put.)
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

# destRow -> sourceRow : destRow should end up with the payload that
# currently (pre-edit) lives on sourceRow.
$map = @{
    14 = 15; 15 = 14;
    23 = 24; 24 = 23;
    33 = 34; 34 = 33;
    36 = 38; 37 = 36; 38 = 37;
    54 = 55; 55 = 54;
    72 = 77; 73 = 74; 74 = 73; 75 = 72; 76 = 75; 77 = 76;
    85 = 86; 86 = 85;
    92 = 93; 93 = 95; 95 = 92;
    102 = 105; 103 = 106; 104 = 102; 105 = 104; 106 = 103;
}

# 1) Snapshot every affected row's payload BEFORE any writes happen -
#    several of the remaps form cycles (e.g. 72->77->76->75->72...),
#    so we must not read a row after it has already been overwritten.
$snapshot = @{}
foreach ($row in $map.Keys) {
    $values = @()
    foreach ($col in $cols) {
        $values += , ($ws.Range("$col$row").Value2)
    }
    $snapshot[$row] = $values
}

# 2) Write each destination row's columns from the snapshot taken from
#    its source row.
foreach ($row in $map.Keys) {
    $srcRow = $map[$row]
    $values = $snapshot[$srcRow]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
